$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.44423134029726
$ws.Range("C2").Value = 8.966419816057178
$ws.Range("E2").Value = 12.13397574667697
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 3.663921409850964
$ws.Range("K2").Value = 9.805396390516554
$ws.Range("L2").Value = 9.680989646931536
$ws.Range("N2").Value = 19.82714434424692
$ws.Range("O2").Value = 25.36111982046866

$ws.Range("B3").Value = 13.22973205327845
$ws.Range("C3").Value = 8.975085404275843
$ws.Range("E3").Value = 12.14133167830885
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 3.66560950315069
$ws.Range("K3").Value = 9.657613798352548
$ws.Range("L3").Value = 9.667816930363271
$ws.Range("N3").Value = 19.88856930268932
$ws.Range("O3").Value = 25.45260791865534

$ws.Range("B4").Value = 13.09884643850672
$ws.Range("C4").Value = 8.980879977038763
$ws.Range("E4").Value = 12.14811255993619
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 3.666701200892106
$ws.Range("K4").Value = 9.567302026380577
$ws.Range("L4").Value = 9.661362693579383
$ws.Range("N4").Value = 19.92808527243913
$ws.Range("O4").Value = 25.51357946385107

$ws.Range("B5").Value = 13.04578376993756
$ws.Range("C5").Value = 8.983360844958897
$ws.Range("E5").Value = 12.15144569349812
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 3.667160000001231
$ws.Range("K5").Value = 9.530650969242799
$ws.Range("L5").Value = 9.659145736080189
$ws.Range("N5").Value = 19.94464249567014
$ws.Range("O5").Value = 25.53963089766154

$ws.Range("B6").Value = 13.03699132938315
$ws.Range("C6").Value = 8.983780021662815
$ws.Range("E6").Value = 12.15203358688532
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 3.667237025435435
$ws.Range("K6").Value = 9.524575548421595
$ws.Range("L6").Value = 9.658802631470692
$ws.Range("N6").Value = 19.94741927804918
$ws.Range("O6").Value = 25.54402946038891

$ws.Range("B7").Value = 13.09812961542908
$ws.Range("C7").Value = 8.980912950407834
$ws.Range("E7").Value = 12.14815520385745
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 3.666707331985452
$ws.Range("K7").Value = 9.566807064527145
$ws.Range("L7").Value = 9.661331118989665
$ws.Range("N7").Value = 19.92830672837815
$ws.Range("O7").Value = 25.51392592523742

$ws.Range("B8").Value = 13.37014333050569
$ws.Range("C8").Value = 8.969309584695427
$ws.Range("E8").Value = 12.13604246517948
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 3.664492032959779
$ws.Range("K8").Value = 9.754378758633685
$ws.Range("L8").Value = 9.676110070191136
$ws.Range("N8").Value = 19.84795071205416
$ws.Range("O8").Value = 25.39166840439112

$ws.Range("B9").Value = 13.90706020814482
$ws.Range("C9").Value = 8.950298858805583
$ws.Range("E9").Value = 12.13022707364345
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.660583896476795
$ws.Range("K9").Value = 10.12366755383578
$ws.Range("L9").Value = 9.717942006504439
$ws.Range("N9").Value = 19.70460021968191
$ws.Range("O9").Value = 25.19004670135918

$ws.Range("B10").Value = 14.29973486799456
$ws.Range("C10").Value = 8.938591603483975
$ws.Range("E10").Value = 12.13684072461828
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.657975700330422
$ws.Range("K10").Value = 10.39333471308695
$ws.Range("L10").Value = 9.756342715726884
$ws.Range("N10").Value = 19.60786723400752
$ws.Range("O10").Value = 25.06523720516713

$ws.Range("B11").Value = 14.47718181049824
$ws.Range("C11").Value = 8.933751698836563
$ws.Range("E11").Value = 12.14219845734109
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.656845718019391
$ws.Range("K11").Value = 10.51513619437518
$ws.Range("L11").Value = 9.775435984429919
$ws.Range("N11").Value = 19.56570690181917
$ws.Range("O11").Value = 25.01353925808382

$ws.Range("B12").Value = 14.54414249121511
$ws.Range("C12").Value = 8.931988427671051
$ws.Range("E12").Value = 12.14456343813958
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.656425903926239
$ws.Range("K12").Value = 10.56109240378893
$ws.Range("L12").Value = 9.782895819969587
$ws.Range("N12").Value = 19.55000567867084
$ws.Range("O12").Value = 24.99469431807761

$ws.Range("B13").Value = 14.52973274473445
$ws.Range("C13").Value = 8.932365094545041
$ws.Range("E13").Value = 12.14403917722617
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.656515959404194
$ws.Range("K13").Value = 10.55120300777826
$ws.Range("L13").Value = 9.781279064201973
$ws.Range("N13").Value = 19.55337549615906
$ws.Range("O13").Value = 24.99872033312767

$ws.Range("B14").Value = 14.48269576843666
$ws.Range("C14").Value = 8.933605242716499
$ws.Range("E14").Value = 12.14238629975818
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.656811017833586
$ws.Range("K14").Value = 10.51892062877226
$ws.Range("L14").Value = 9.776045131192824
$ws.Range("N14").Value = 19.56440986893266
$ws.Range("O14").Value = 25.01197419469481

$ws.Range("B15").Value = 14.45385183702013
$ws.Range("C15").Value = 8.934373909011098
$ws.Range("E15").Value = 12.14141758400796
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.656992801488423
$ws.Range("K15").Value = 10.49912372284374
$ws.Range("L15").Value = 9.772868978045167
$ws.Range("N15").Value = 19.57120308388792
$ws.Range("O15").Value = 25.0201879381905

$ws.Range("B16").Value = 14.28810891527522
$ws.Range("C16").Value = 8.938917647110713
$ws.Range("E16").Value = 12.13653771660309
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.658050680913246
$ws.Range("K16").Value = 10.38535353246765
$ws.Range("L16").Value = 9.75512728171808
$ws.Range("N16").Value = 19.61065951517022
$ws.Range("O16").Value = 25.06871812331263

$ws.Range("B17").Value = 14.18608132042714
$ws.Range("C17").Value = 8.941829248275713
$ws.Range("E17").Value = 12.1341446261789
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 3.658714098742006
$ws.Range("K17").Value = 10.31530558105349
$ws.Range("L17").Value = 9.744656667251412
$ws.Range("N17").Value = 19.63533620139749
$ws.Range("O17").Value = 25.09979174319633

$ws.Range("B18").Value = 14.12729012526843
$ws.Range("C18").Value = 8.943549670055527
$ws.Range("E18").Value = 12.13298941208743
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 3.659100999402662
$ws.Range("K18").Value = 10.27493624577866
$ws.Range("L18").Value = 9.73878741584541
$ws.Range("N18").Value = 19.64970321350471
$ws.Range("O18").Value = 25.11814238347641

$ws.Range("B19").Value = 14.10736796253965
$ws.Range("C19").Value = 8.944140043158397
$ws.Range("E19").Value = 12.13263632007497
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 3.659232912269347
$ws.Range("K19").Value = 10.26125551743118
$ws.Range("L19").Value = 9.736826612427848
$ws.Range("N19").Value = 19.65459749517295
$ws.Range("O19").Value = 25.12443762637491

$ws.Range("B20").Value = 14.19695394644563
$ws.Range("C20").Value = 8.941514571453579
$ws.Range("E20").Value = 12.13437648990539
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 3.658642926435427
$ws.Range("K20").Value = 10.3227708537319
$ws.Range("L20").Value = 9.745755457687784
$ws.Range("N20").Value = 19.63269136406951
$ws.Range("O20").Value = 25.09643442864646

$ws.Range("B21").Value = 14.49651853735697
$ws.Range("C21").Value = 8.933239097906021
$ws.Range("E21").Value = 12.14286268217031
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.656724132885305
$ws.Range("K21").Value = 10.52840761540905
$ws.Range("L21").Value = 9.777576265118981
$ws.Range("N21").Value = 19.56116165121279
$ws.Range("O21").Value = 25.00806133619305

$ws.Range("B22").Value = 14.69090573880504
$ws.Range("C22").Value = 8.928235505260952
$ws.Range("E22").Value = 12.1503671440467
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.655517202145496
$ws.Range("K22").Value = 10.66180925090377
$ws.Range("L22").Value = 9.799709433895016
$ws.Range("N22").Value = 19.51595098360948
$ws.Range("O22").Value = 24.95457153449928

$ws.Range("B23").Value = 14.58730594753765
$ws.Range("C23").Value = 8.930869086338991
$ws.Range("E23").Value = 12.14618330653797
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 3.656157065810643
$ws.Range("K23").Value = 10.59071469233229
$ws.Range("L23").Value = 9.787775642052475
$ws.Range("N23").Value = 19.53994041629375
$ws.Range("O23").Value = 24.9827290624802

$ws.Range("B24").Value = 14.19203884518529
$ws.Range("C24").Value = 8.941656691969897
$ws.Range("E24").Value = 12.13427097703652
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.658675086322128
$ws.Range("K24").Value = 10.31939610546063
$ws.Range("L24").Value = 9.745258225619157
$ws.Range("N24").Value = 19.63388653405962
$ws.Range("O24").Value = 25.09795075680209

$ws.Range("B25").Value = 13.76184708001262
$ws.Range("C25").Value = 8.955043297960549
$ws.Range("E25").Value = 12.12988413952323
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 3.661594750400462
$ws.Range("K25").Value = 10.02387443531649
$ws.Range("L25").Value = 9.705265833611969
$ws.Range("N25").Value = 19.74186619972813
$ws.Range("O25").Value = 25.24050086265823
